# Pull in the latest columns from the DB export: a new "04_05_20" date
# column (D) alongside the existing "03_05_20" column (C), plus the
# matching data value for the existing student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added date column
$ws.Range("D1").Value = "04_05_20"

# Data point for the existing row (numeric, unlike the inline-string cells)
$ws.Range("D3").Value = 1

# The blank spacer row's cells were placeholder inline strings with no
# text; clear them out now that the sheet has real shared-string content.
$ws.Range("A2:C2").ClearContents()
